# Apply "Default Color Scheme" palette swap:
#   - Slide master background: navy (000080) -> brick red (A13942)
#   - All title/body run text color: sienna (A0522D) -> amber/gold (BD8038)

$p = $ppt.ActivePresentation

# 1. Update the slide master background fill color.
$master = $p.SlideMaster
$master.Background.Fill.ForeColor.RGB = 4340129   # 0xA13942

# 2. Update every shape's text color across all slides.
for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $slide = $p.Slides.Item($i)
    for ($j = 1; $j -le $slide.Shapes.Count; $j++) {
        $shape = $slide.Shapes.Item($j)
        if ($shape.HasTextFrame) {
            $tr = $shape.TextFrame.TextRange
            if ($tr.Font.Color.RGB -eq 2970272) {   # 0xA0522D
                $tr.Font.Color.RGB = 3702973        # 0xBD8038
            }
        }
    }
}
